$wb = $excel.ActiveWorkbook
Write-Host "Sheets:"
foreach ($ws in $wb.Worksheets) {
    Write-Host $ws.Name
}
